$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.173.53'
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").Value = '2.560.40'
$ws.Range("E3").Value = '  +0.22%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = "'583.62"
$ws.Range("E5").Value = '  +2.86%  '
$ws.Range("D6").Value = "'147.50"
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = "'0.588"
$ws.Range("E8").Value = '  +1.43%  '
$ws.Range("E9").Value = '  +3.97%  '
$ws.Range("E10").Value = '  +0.94%  '
$ws.Range("E11").Value = '  +0.40%  '
$ws.Range("E12").Value = '  +0.73%  '
$ws.Range("D13").Value = "'27.55"
$ws.Range("E13").Value = '  +1.61%  '
$ws.Range("D14").Value = '3.019.27'
$ws.Range("E14").Value = '  +0.18%  '
$ws.Range("D15").Value = '63.105.35'
$ws.Range("E15").Value = '  +0.31%  '
$ws.Range("E16").Value = '  +5.21%  '
$ws.Range("D17").Value = '2.548.62'
$ws.Range("E17").Value = '  +0.96%  '
$ws.Range("D18").Value = "'11.34"
$ws.Range("E18").Value = '  -1.23%  '
$ws.Range("E19").Value = '  +3.83%  '
$ws.Range("D20").Value = "'341.42"
$ws.Range("E20").Value = '  +2.17%  '
$ws.Range("D21").Value = "'6.83"
$ws.Range("E21").Value = '  +1.03%  '
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("D23").Value = "'66.58"
$ws.Range("E23").Value = '  +2.50%  '
$ws.Range("B24").Value = 'WrappedeETH'
$ws.Range("C24").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D24").Value = '2.684.17'
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("B25").Value = 'Fetch.AI'
$ws.Range("C25").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D25").Value = "'1.64"
$ws.Range("E25").Value = '  +2.82%  '
$ws.Range("D26").Value = "'0.170"
$ws.Range("E26").Value = '  +1.03%  '
$ws.Range("D27").Value = "'8.16"
$ws.Range("E27").Value = '  +12.61%  '
$ws.Range("D28").Value = "'8.56"
$ws.Range("E28").Value = '  +2.28%  '
$ws.Range("B29").Value = 'SuiNetwork'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D29").Value = "'1.49"
$ws.Range("E29").Value = '  +0.50%  '
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").Value = "'0.996"
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("D31").Value = "'2.02"
$ws.Range("E31").Value = '  +8.94%  '
$ws.Range("E32").Value = '  +1.91%  '
$ws.Range("D33").Value = "'457.21"
$ws.Range("E33").Value = '  +11.40%  '
$ws.Range("D34").Value = "'176.78"
$ws.Range("E34").Value = '  -0.09%  '
$ws.Range("D35").Value = "'1.62"
$ws.Range("E35").Value = '  +3.08%  '
$ws.Range("E36").Value = '  +2.53%  '
$ws.Range("D37").Value = "'19.30"
$ws.Range("E37").Value = '  +2.64%  '
$ws.Range("D38").Value = "'4.53"
$ws.Range("E38").Value = '  +3.87%  '
$ws.Range("E39").Value = '  +0.05%  '
$ws.Range("D40").Value = "'1.74"
$ws.Range("E40").Value = '  +0.25%  '
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").Value = "'151.08"
$ws.Range("E42").Value = '  -0.47%  '
$ws.Range("D43").Value = "'3.83"
$ws.Range("E43").Value = '  +2.44%  '
$ws.Range("D44").Value = "'21.10"
$ws.Range("E44").Value = '  +2.25%  '
$ws.Range("E45").Value = '  +6.88%  '
$ws.Range("D46").Value = "'0.616"
$ws.Range("E46").Value = '  +1.76%  '
$ws.Range("D47").Value = "'0.0978"
$ws.Range("E47").Value = '  +2.13%  '
$ws.Range("D48").Value = "'0.0241"
$ws.Range("E48").Value = '  +2.27%  '
$ws.Range("D49").Value = "'18.44"
$ws.Range("E49").Value = '  +0.25%  '
$ws.Range("E50").Value = '  -1.55%  '
$ws.Range("D51").Value = "'11.39"
$ws.Range("E51").Value = '  -0.20%  '
